# Week 17 data: a new player column ("S.Williams") is logged on both the
# "Rushing" and "Receiving" sheets. It is inserted as a new column just
# before the existing "N.Fant" column (column Q), pushing N.Fant and
# everything to its right one column to the right (through a new column U).
# The new column gets the same header styling as the other player columns
# and an "n" placeholder in the data row, matching the rest of that row.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Rushing", "Receiving")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new column at Q, shifting N.Fant (and everything after it)
    # one column to the right.
    $ws.Columns("Q:Q").Insert()

    # Populate the newly inserted column's header + data cells.
    $ws.Range("Q1").Value = "S.Williams"
    $ws.Range("Q2").Value = "n"
}
